$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Enter the new "Tiempo Real" (T12) value for row 12; dependent shared
# formulas (U12, X12, AA12, ... BA12) recalculate automatically.
$ws.Range("T12").Value = 5

# Re-merge the header cells so the mergeCells list ends up with the
# AZ4:BA4 / AO4:AP4 / AR4:AS4 / AU4:AV4 / AX4:AY4 ranges first, matching
# the sheet as re-serialized after the merged header cells were last
# touched. Re-merging (unmerge then merge) a range moves its entry to
# the end of the mergeCells list, so re-touching the eleven ranges
# that must end up *after* the first five is the minimal sequence that
# reproduces the target order.
$ws.Range("AL4:AM4").UnMerge()
$ws.Range("H4:I4").UnMerge()
$ws.Range("K4:L4").UnMerge()
$ws.Range("N4:O4").UnMerge()
$ws.Range("Q4:R4").UnMerge()
$ws.Range("T4:U4").UnMerge()
$ws.Range("W4:X4").UnMerge()
$ws.Range("Z4:AA4").UnMerge()
$ws.Range("AC4:AD4").UnMerge()
$ws.Range("AF4:AG4").UnMerge()
$ws.Range("AI4:AJ4").UnMerge()

$ws.Range("AL4:AM4").Merge()
$ws.Range("H4:I4").Merge()
$ws.Range("K4:L4").Merge()
$ws.Range("N4:O4").Merge()
$ws.Range("Q4:R4").Merge()
$ws.Range("T4:U4").Merge()
$ws.Range("W4:X4").Merge()
$ws.Range("Z4:AA4").Merge()
$ws.Range("AC4:AD4").Merge()
$ws.Range("AF4:AG4").Merge()
$ws.Range("AI4:AJ4").Merge()

# Update the active selection in the bottom-right frozen pane.
$ws.Range("F12").Select()
